$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.028.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.585.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.77%  '

$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.004'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3764'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3582'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.04%  '

$ws.Range("E10").Value = '  -0.13%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08047'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.216'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.93'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.07%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.461'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.306'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.34%  '

$ws.Range("E16").Value = '  -2.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.588.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.97%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06801'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.429'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.018.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.368'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.772'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.44%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.193'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.65%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.351'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.520'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.766.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9383'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07347'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02673'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.977'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08759'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.058'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.97%  '

$ws.Range("E40").Value = '  -2.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.336'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6875'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.84'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6401'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.988'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.80%  '

$ws.Range("E47").Value = '  -2.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07883'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.190'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.190'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.71%  '
